$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2122448979591837
$ws.Range("C2").Value = 0.5265306122448979
$ws.Range("J2").Value = 0.00816326530612245
$ws.Range("P2").Value = 0.1551020408163265
$ws.Range("S2").Value = 0.09795918367346938
$ws.Range("B3").Value = 0.007352941176470588
$ws.Range("C3").Value = 0.02205882352941177
$ws.Range("J3").Value = 0.01470588235294118
$ws.Range("P3").Value = 0.6397058823529411
$ws.Range("S3").Value = 0.3161764705882353
$ws.Range("B6").Value = 0.02475247524752475
$ws.Range("D6").Value = 0.0297029702970297
$ws.Range("F6").Value = 0.0594059405940594
$ws.Range("J6").Value = 0.2524752475247525
$ws.Range("O6").Value = 0.0297029702970297
$ws.Range("Q6").Value = 0.1386138613861386
$ws.Range("R6").Value = 0.08415841584158416
$ws.Range("S6").Value = 0.3811881188118812
$ws.Range("B7").Value = 0.1173184357541899
$ws.Range("D7").Value = 0.0111731843575419
$ws.Range("F7").Value = 0.05027932960893855
$ws.Range("J7").Value = 0.1284916201117318
$ws.Range("O7").Value = 0.01675977653631285
$ws.Range("Q7").Value = 0.1452513966480447
$ws.Range("R7").Value = 0.0782122905027933
$ws.Range("S7").Value = 0.4525139664804469
$ws.Range("B8").Value = 0.08983451536643026
$ws.Range("D8").Value = 0.01182033096926714
$ws.Range("E8").Value = 0.002364066193853428
$ws.Range("F8").Value = 0.04728132387706856
$ws.Range("J8").Value = 0.1205673758865248
$ws.Range("O8").Value = 0.009456264775413711
$ws.Range("Q8").Value = 0.1678486997635934
$ws.Range("R8").Value = 0.09692671394799054
$ws.Range("S8").Value = 0.4539007092198581
$ws.Range("B9").Value = 0.08500000000000001
$ws.Range("D9").Value = 0.015
$ws.Range("F9").Value = 0.06
$ws.Range("J9").Value = 0.12
$ws.Range("O9").Value = 0.025
$ws.Range("Q9").Value = 0.14
$ws.Range("R9").Value = 0.08
$ws.Range("S9").Value = 0.475
$ws.Range("B10").Value = 0.08930075821398484
$ws.Range("D10").Value = 0.01769165964616681
$ws.Range("E10").Value = 0.002527379949452401
$ws.Range("F10").Value = 0.07413647851727043
$ws.Range("J10").Value = 0.1196293176074136
$ws.Range("O10").Value = 0.01769165964616681
$ws.Range("Q10").Value = 0.1996630160067397
$ws.Range("R10").Value = 0.09856781802864364
$ws.Range("S10").Value = 0.3807919123841618
$ws.Range("G11").Value = 0.1254237288135593
$ws.Range("J11").Value = 0.1016949152542373
$ws.Range("K11").Value = 0.1830508474576271
$ws.Range("L11").Value = 0.5694915254237288
$ws.Range("S11").Value = 0.02033898305084746
$ws.Range("G12").Value = 0.7241379310344828
$ws.Range("J12").Value = 0.2126436781609195
$ws.Range("K12").Value = 0.005747126436781609
$ws.Range("L12").Value = 0.02298850574712644
$ws.Range("S12").Value = 0.03448275862068965
$ws.Range("G13").Value = 0.5405405405405406
$ws.Range("J13").Value = 0.3783783783783784
$ws.Range("S13").Value = 0.08108108108108109
$ws.Range("F15").Value = 0.0108695652173913
$ws.Range("H15").Value = 0.1304347826086956
$ws.Range("I15").Value = 0.07065217391304347
$ws.Range("J15").Value = 0.3315217391304348
$ws.Range("K15").Value = 0.05434782608695652
$ws.Range("M15").Value = 0.01630434782608696
$ws.Range("O15").Value = 0.07065217391304347
$ws.Range("S15").Value = 0.3152173913043478
$ws.Range("F16").Value = 0.02054794520547945
$ws.Range("H16").Value = 0.1438356164383562
$ws.Range("I16").Value = 0.08904109589041095
$ws.Range("J16").Value = 0.4246575342465753
$ws.Range("K16").Value = 0.1301369863013699
$ws.Range("M16").Value = 0.00684931506849315
$ws.Range("O16").Value = 0.03424657534246575
$ws.Range("S16").Value = 0.1506849315068493
$ws.Range("F17").Value = 0.01288659793814433
$ws.Range("H17").Value = 0.134020618556701
$ws.Range("I17").Value = 0.09020618556701031
$ws.Range("J17").Value = 0.4561855670103093
$ws.Range("K17").Value = 0.09793814432989691
$ws.Range("M17").Value = 0.01288659793814433
$ws.Range("N17").Value = 0.002577319587628866
$ws.Range("O17").Value = 0.05927835051546392
$ws.Range("S17").Value = 0.134020618556701
$ws.Range("F18").Value = 0.01470588235294118
$ws.Range("H18").Value = 0.1813725490196078
$ws.Range("I18").Value = 0.08823529411764706
$ws.Range("J18").Value = 0.3529411764705883
$ws.Range("K18").Value = 0.1323529411764706
$ws.Range("M18").Value = 0.009803921568627451
$ws.Range("O18").Value = 0.06862745098039216
$ws.Range("S18").Value = 0.1519607843137255
$ws.Range("F19").Value = 0.01472868217054264
$ws.Range("H19").Value = 0.2271317829457364
$ws.Range("I19").Value = 0.09612403100775194
$ws.Range("J19").Value = 0.3527131782945737
$ws.Range("K19").Value = 0.1108527131782946
$ws.Range("M19").Value = 0.02170542635658915
$ws.Range("N19").Value = 0.001550387596899225
$ws.Range("O19").Value = 0.05503875968992248
$ws.Range("S19").Value = 0.1201550387596899
